# Update exam time slots on the active worksheet.
# For each exam date (group of consecutive rows sharing the same date in
# column A), the start/end times (columns B/C) are re-sequenced into
# back-to-back slots: 09:00-10:15, 10:15-11:30, 11:45-13:00, 13:15-14:30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$slots = @(
    @("09:00", "10:15"),
    @("10:15", "11:30"),
    @("11:45", "13:00"),
    @("13:15", "14:30")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$slotIndex = 0
$prevDate = $null

for ($r = 2; $r -le $lastRow; $r++) {
    $currentDate = $ws.Cells.Item($r, 1).Value()

    if ($prevDate -ne $null -and $currentDate -ne $prevDate) {
        $slotIndex = 0
    }

    $slot = $slots[$slotIndex]
    $ws.Cells.Item($r, 2).Value = $slot[0]
    $ws.Cells.Item($r, 3).Value = $slot[1]

    $slotIndex++
    $prevDate = $currentDate
}
